{"js": "// Daily-practice sheet refresh: update the date label and the 25\n// two-digit-by-two-digit multiplication problems/answers in the table.\n// Each entry is [oldText, newText]; every value in the document is\n// unique, so an exact-text search reliably targets the right run.\nconst replacements = [\n  [\"2024-12-22 Sunday\", \"2024-12-23 Monday\"],\n  [\"97\u00d786=8342\", \"74\u00d758=4292\"],\n  [\"29\u00d794=2726\", \"18\u00d747=846\"],\n  [\"31\u00d739=1209\", \"78\u00d764=4992\"],\n  [\"12\u00d763=756\", \"68\u00d730=2040\"],\n  [\"19\u00d717=323\", \"94\u00d729=2726\"],\n  [\"73\u00d748=3504\", \"21\u00d750=1050\"],\n  [\"16\u00d787=1392\", \"73\u00d787=6351\"],\n  [\"91\u00d738=3458\", \"33\u00d749=1617\"],\n  [\"16\u00d788=1408\", \"66\u00d791=6006\"],\n  [\"32\u00d743=1376\", \"62\u00d715=930\"],\n  [\"47\u00d745=2115\", \"15\u00d756=840\"],\n  [\"14\u00d775=1050\", \"15\u00d769=1035\"],\n  [\"56\u00d797=5432\", \"90\u00d799=8910\"],\n  [\"38\u00d724=912\", \"56\u00d765=3640\"],\n  [\"77\u00d784=6468\", \"81\u00d711=891\"],\n  [\"20\u00d735=700\", \"58\u00d763=3654\"],\n  [\"84\u00d769=5796\", \"52\u00d741=2132\"],\n  [\"54\u00d795=5130\", \"55\u00d738=2090\"],\n  [\"23\u00d747=1081\", \"13\u00d714=182\"],\n  [\"57\u00d785=4845\", \"40\u00d727=1080\"],\n  [\"62\u00d721=1302\", \"83\u00d762=5146\"],\n  [\"34\u00d752=1768\", \"78\u00d783=6474\"],\n  [\"63\u00d713=819\", \"71\u00d746=3266\"],\n  [\"73\u00d731=2263\", \"55\u00d715=825\"],\n  [\"62\u00d734=2108\", \"51\u00d769=3519\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + oldText);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and the 25 two-digit multiplication answers\n# in the document body, matching the exact before/after text pairs.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-12-22 Sunday\", \"2024-12-23 Monday\"),\n  @(\"97\u00d786=8342\", \"74\u00d758=4292\"),\n  @(\"29\u00d794=2726\", \"18\u00d747=846\"),\n  @(\"31\u00d739=1209\", \"78\u00d764=4992\"),\n  @(\"12\u00d763=756\", \"68\u00d730=2040\"),\n  @(\"19\u00d717=323\", \"94\u00d729=2726\"),\n  @(\"73\u00d748=3504\", \"21\u00d750=1050\"),\n  @(\"16\u00d787=1392\", \"73\u00d787=6351\"),\n  @(\"91\u00d738=3458\", \"33\u00d749=1617\"),\n  @(\"16\u00d788=1408\", \"66\u00d791=6006\"),\n  @(\"32\u00d743=1376\", \"62\u00d715=930\"),\n  @(\"47\u00d745=2115\", \"15\u00d756=840\"),\n  @(\"14\u00d775=1050\", \"15\u00d769=1035\"),\n  @(\"56\u00d797=5432\", \"90\u00d799=8910\"),\n  @(\"38\u00d724=912\", \"56\u00d765=3640\"),\n  @(\"77\u00d784=6468\", \"81\u00d711=891\"),\n  @(\"20\u00d735=700\", \"58\u00d763=3654\"),\n  @(\"84\u00d769=5796\", \"52\u00d741=2132\"),\n  @(\"54\u00d795=5130\", \"55\u00d738=2090\"),\n  @(\"23\u00d747=1081\", \"13\u00d714=182\"),\n  @(\"57\u00d785=4845\", \"40\u00d727=1080\"),\n  @(\"62\u00d721=1302\", \"83\u00d762=5146\"),\n  @(\"34\u00d752=1768\", \"78\u00d783=6474\"),\n  @(\"63\u00d713=819\", \"71\u00d746=3266\"),\n  @(\"73\u00d731=2263\", \"55\u00d715=825\"),\n  @(\"62\u00d734=2108\", \"51\u00d769=3519\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n  #          MatchAllWordForms, Forward, Wrap(=wdFindContinue), Format, ReplaceWith,\n  #          Replace(=wdReplaceAll))\n  $ok = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $ok) {\n    throw \"Could not find expected text: $oldText\"\n  }\n}\n"}
